$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G3").Value = 3.4
$ws.Range("I3").Value = 2.5
$ws.Range("G6").Value = 1.7
$ws.Range("H6").Value = 3.9
$ws.Range("I6").Value = 4.75
$ws.Range("R6").Value = 1.75
$ws.Range("S6").Value = 2
$ws.Range("U6").Value = 8.5
$ws.Range("W6").Value = 13
$ws.Range("AA6").Value = 7.5
$ws.Range("AB6").Value = 15
$ws.Range("AD6").Value = 201
$ws.Range("AJ6").Value = 41
$ws.Range("G7").Value = 1.25
$ws.Range("H7").Value = 4.8
$ws.Range("I7").Value = 11.75
$ws.Range("L7").Value = 1.22
$ws.Range("M7").Value = 3.45
$ws.Range("N7").Value = 1.65
$ws.Range("O7").Value = 1.98
$ws.Range("R7").Value = 2.18
$ws.Range("S7").Value = 1.53
$ws.Range("T7").Value = 6.2
$ws.Range("U7").Value = 5.5
$ws.Range("V7").Value = 9
$ws.Range("W7").Value = 7.1
$ws.Range("X7").Value = 11.5
$ws.Range("Z7").Value = 11.5
$ws.Range("AA7").Value = 10.25
$ws.Range("AB7").Value = 26
$ws.Range("AE7").Value = 28
$ws.Range("AF7").Value = 100
$ws.Range("AG7").Value = 37
$ws.Range("AH7").Value = 450
$ws.Range("AI7").Value = 175
$ws.Range("G8").Value = 1.32
$ws.Range("H8").Value = 4.55
$ws.Range("I8").Value = 8.75
$ws.Range("M8").Value = 3.15
$ws.Range("N8").Value = 1.78
$ws.Range("O8").Value = 1.82
$ws.Range("R8").Value = 2.2
$ws.Range("W8").Value = 7.6
$ws.Range("X8").Value = 12.5
$ws.Range("Z8").Value = 10.25
$ws.Range("AA8").Value = 9.25
$ws.Range("AC8").Value = 150
$ws.Range("AE8").Value = 19.5
$ws.Range("AF8").Value = 60
$ws.Range("AG8").Value = 29
$ws.Range("AH8").Value = 250
$ws.Range("AI8").Value = 120
$ws.Range("G12").Value = 2.92
$ws.Range("H12").Value = 2.37
$ws.Range("J12").Value = 1.19
$ws.Range("K12").Value = 4.1
$ws.Range("L12").Value = 1.72
$ws.Range("M12").Value = 2
$ws.Range("N12").Value = 3.1
$ws.Range("O12").Value = 1.32
$ws.Range("P12").Value = 1.72
$ws.Range("Q12").Value = 2
$ws.Range("R12").Value = 2.25
$ws.Range("V12").Value = 11.5
$ws.Range("Z12").Value = 4.1
$ws.Range("AB12").Value = 18.5
$ws.Range("AE12").Value = 6.2
$ws.Range("H13").Value = 2.77
$ws.Range("K13").Value = 4.55
$ws.Range("L13").Value = 1.7
$ws.Range("N13").Value = 3.05
$ws.Range("O13").Value = 1.33
$ws.Range("T13").Value = 5
$ws.Range("V13").Value = 10.5
$ws.Range("X13").Value = 27
$ws.Range("Y13").Value = 60
$ws.Range("Z13").Value = 4.55
$ws.Range("AE13").Value = 6.7
$ws.Range("AG13").Value = 14.5
$ws.Range("I14").Value = 3.1
$ws.Range("N14").Value = 2.5
$ws.Range("O14").Value = 1.5
$ws.Range("U14").Value = 10
$ws.Range("AE14").Value = 7.5
$ws.Range("G16").Value = 2.15
$ws.Range("I16").Value = 3.3
$ws.Range("N16").Value = 2.08
$ws.Range("O16").Value = 1.73
$ws.Range("U16").Value = 10
$ws.Range("AH16").Value = 34
$ws.Range("I17").Value = 4
$ws.Range("N17").Value = 2.05
$ws.Range("O17").Value = 1.75
$ws.Range("P17").Value = 1.4
$ws.Range("Q17").Value = 2.75
$ws.Range("R17").Value = 1.91
$ws.Range("S17").Value = 1.91
$ws.Range("T17").Value = 7
$ws.Range("V17").Value = 8.5
$ws.Range("Z17").Value = 9
$ws.Range("AD17").Value = 301
$ws.Range("AE17").Value = 11
$ws.Range("AG17").Value = 13
$ws.Range("AI17").Value = 34
$ws.Range("G18").Value = 2.1
$ws.Range("I18").Value = 3.5
$ws.Range("L18").Value = 1.44
$ws.Range("M18").Value = 2.63
$ws.Range("T18").Value = 6
$ws.Range("X18").Value = 19
$ws.Range("G19").Value = 2.5
$ws.Range("I19").Value = 3.2
$ws.Range("U19").Value = 10
$ws.Range("G20").Value = 1.8
$ws.Range("H20").Value = 3.5
$ws.Range("I20").Value = 4.5
$ws.Range("J20").Value = 1.07
$ws.Range("K20").Value = 9
$ws.Range("P20").Value = 1.4
$ws.Range("Q20").Value = 2.75
$ws.Range("V20").Value = 8.5
$ws.Range("X20").Value = 15
$ws.Range("Z20").Value = 9
$ws.Range("AH20").Value = 51
$ws.Range("L21").Value = 1.3
$ws.Range("M21").Value = 3.4
$ws.Range("N21").Value = 2
$ws.Range("O21").Value = 1.8
$ws.Range("N22").Value = 1.85
$ws.Range("O22").Value = 1.95
$ws.Range("I24").Value = 2.7
$ws.Range("R24").Value = 2
$ws.Range("S24").Value = 1.73
$ws.Range("W24").Value = 29
$ws.Range("X24").Value = 26
$ws.Range("Z24").Value = 7
$ws.Range("AB24").Value = 17
$ws.Range("G25").Value = 2.3
$ws.Range("I25").Value = 3.5
$ws.Range("J25").Value = 1.1
$ws.Range("K25").Value = 7
$ws.Range("V25").Value = 10
$ws.Range("AH25").Value = 34
$ws.Range("AI25").Value = 29
$ws.Range("J26").Value = 1.05
$ws.Range("K26").Value = 11
$ws.Range("L26").Value = 1.3
$ws.Range("M26").Value = 3.4
$ws.Range("N26").Value = 1.98
$ws.Range("O26").Value = 1.83
$ws.Range("T26").Value = 9.5
$ws.Range("U26").Value = 17
$ws.Range("V26").Value = 12
$ws.Range("AB26").Value = 15
$ws.Range("AI26").Value = 17
$ws.Range("G29").Value = 2.15
$ws.Range("I29").Value = 3.2
$ws.Range("P29").Value = 1.36
$ws.Range("Q29").Value = 3
$ws.Range("R29").Value = 1.73
$ws.Range("S29").Value = 2
$ws.Range("T29").Value = 8
$ws.Range("Y29").Value = 26
$ws.Range("Z29").Value = 11
$ws.Range("AD29").Value = 201
$ws.Range("AE29").Value = 10
$ws.Range("AG29").Value = 12
$ws.Range("AI29").Value = 26
$ws.Range("AJ29").Value = 34
$ws.Range("G42").Value = 2.2
$ws.Range("K42").Value = 12
$ws.Range("N42").Value = 1.8
$ws.Range("O42").Value = 2
